$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("What's the maximum number of lithology types in an log?", "llama3.2:latest", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("What's the maximum number of lithology types in an log?", "deepseek1.5", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("How many tracks can you define in one ODF?", "deepseek1.5", "According to the Track Settings dialog box, the number of tracks that can be defined is 200."),
    @("How many curve shades can I create?", "deepseek1.5", "According to the document, you can create 250 curve shades."),
    @("What's the maximum number of lithology types in an log?", "llama3.2:latest", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("How many tracks can you define in one ODF?", "llama3.2:latest", "According to the Track Settings dialog box, the number of tracks that can be defined is 200."),
    @("How many curve shades can I create?", "llama3.2:latest", "According to the document, you can create 250 curve shades."),
    @("How many tables can I have in my log?", "llama3.2:latest", "You can have up to 100 tables in a log.")
)

$startRow = 80
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
